$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stale "x"/"o" capture markers that no longer apply now that
# there are no moves left (keep the cell's existing style/formatting).
$clearCells = @("AC23","AE23","AG23","AI23","AB24","AD24","AF24","AI25","AD26","AD28","AH28","AI29","AB30","AF30","AH30")
foreach ($addr in $clearCells) {
    $ws.Range($addr).Value = ""
}

# These two markers flip from "o" to "x".
$ws.Range("AC27").Value = "x"
$ws.Range("AI27").Value = "x"

# New informational cell showing whose move it is.
$ws.Range("AJ30").Value = "Move: x"

# Update the view so the newly-relevant area is visible/selected.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("AJ31").Select()
